# Add a new school entry (MGA School Kenya) as row 22 on Sheet1, mirroring
# the structure of the existing rows (SchoolName, ShortCode, UserName,
# Password, GradeLabel, Class_Sections, Raw_SectionMap).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "MGA School Kenya"
$ws.Range("B22").Value = "MGA Kenya"
$ws.Range("C22").Value = "kenyamoe1"
$ws.Range("D22").Value = "kenyamoe1"
$ws.Range("E22").Value = "Grade"
$ws.Range("F22").Value = "{01A0=Grade 1 , 02A0=Grade 2, 03A0=Grade 3, 04A0=Grade 4, 05A0=Grade 5, 06A0=Grade 6, 07A0=Grade 7, 08A0=Grade 8, 09A0=Grade 9, 10A0=Grade 10, 11A0=Grade 11, 12A0=Grade 12,}"

# Match the saved view state: the sheet was left scrolled down with G25
# selected (empty cell just past the new row).
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
[void]$ws.Range("G25").Select()
